$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks first; they will be rebuilt after the column shift below
# (the underlying engine does not re-target hyperlink refs on column insert)
$ws.Hyperlinks.Delete()

# Insert a new column before column A.
# This shifts the existing Email column (A) to B, and Password column (B) to C.
$ws.Columns.Item(1).Insert()

# New 7th row of data (email + password entered first, name filled in afterwards)
$ws.Range("B7").Value = "bahsa.leb@gmail.com"
$ws.Range("C7").Value = "1q2w3e4r5t!@#$%"

# Fill in the new "Name" column for the existing rows
$ws.Range("A1").Value = "Samira Raad"
$ws.Range("A2").Value = "Samir Hanna"
$ws.Range("A3").Value = "Marcelle Hanna"
$ws.Range("A4").Value = "Didi Hanna"
$ws.Range("A6").Value = "Marcelle Haddad"
$ws.Range("A5").Value = "Michael Merri"
$ws.Range("A7").Value = "Bahsa Lebron"

# Re-create the hyperlinks, now targeting the shifted email column (B)
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:samira_raad2000@outlook.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:samir_hanna2000@outlook.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:marcelle_hanna@outlook.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:didi_hanna@outlook.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:michael.me2@hotmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:marcelle.haddad@outlook.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:bahsa.leb@gmail.com") | Out-Null

# Adding hyperlinks via COM creates a brand-new cell style; restore the
# shared "Hyperlink" cell style that the other cells already use.
$ws.Range("B1:B7").Style = "Hyperlink"

# Resize columns to fit the new layout (Name / Email / Password)
$ws.Columns.Item(1).ColumnWidth = 17.26
$ws.Columns.Item(2).ColumnWidth = 30.166666666666668
$ws.Columns.Item(3).ColumnWidth = 24.42

# Leave the cursor where the user would naturally end up after entering the data
$ws.Range("C8").Select() | Out-Null
